# Generate Report for handback
# Refresh the "Latest Handoff Datetime" / "Latest Handback DateTime" values
# for the 7b66201e-... source file row (row 2) on both the zh-cn and de-de
# localization-status sheets, reflecting a newly generated handback report.

$wb = $excel.ActiveWorkbook

$ws_zhcn = $wb.Worksheets.Item("zh-cn")
$ws_zhcn.Range("D2").Value = "2016-01-11 03:04:56"
$ws_zhcn.Range("G2").Value = "2016-01-11 03:05:42"

$ws_dede = $wb.Worksheets.Item("de-de")
$ws_dede.Range("D2").Value = "2016-01-11 03:05:07"
$ws_dede.Range("G2").Value = "2016-01-11 03:06:02"
